$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 442; this shifts the existing rows
# 442:513 down to 443:514 and extends the used range to A1:T514.
$ws.Rows("442:442").Insert()

# Populate the newly inserted row 442 with the new record.
$ws.Range("A442").Value = 9
$ws.Range("B442").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C442").Value = "Metropolitana"
$ws.Range("D442").Value = 44522
$ws.Range("E442").Value = 13
$ws.Range("F442").Value = "Fruta"
$ws.Range("G442").Value = 100102
$ws.Range("H442").Value = "Cítricos"
$ws.Range("I442").Value = 100102005
$ws.Range("J442").Value = "Naranja"
$ws.Range("K442").Value = "Valencia"
$ws.Range("L442").Value = "Primera"
$ws.Range("M442").Value = 470
$ws.Range("N442").Value = 8000
$ws.Range("O442").Value = 8500
$ws.Range("P442").Value = 8266
$ws.Range("Q442").Value = "$/malla 18 kilos"
$ws.Range("R442").Value = "Provincia de Quillota"
$ws.Range("S442").Value = 459
$ws.Range("T442").Value = 18
